$d = $word.ActiveDocument

# Insert four new paragraphs at the very beginning of the document body:
#   "Hylde"
#   "*kilder"
#   "https://www.webmatematik.dk/lektioner/7-9-klasse/rumfang-og-overfladeareal/kasse"
#   (empty paragraph)
# followed by the document's existing content (starting with "Vaskemaskine kilder").

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<w:p ' + $wNs + '><w:r><w:t>Hylde</w:t></w:r></w:p>' +
       '<w:p ' + $wNs + '><w:r><w:t>*kilder</w:t></w:r></w:p>' +
       '<w:p ' + $wNs + '><w:r><w:t>https://www.webmatematik.dk/lektioner/7-9-klasse/rumfang-og-overfladeareal/kasse</w:t></w:r></w:p>' +
       '<w:p ' + $wNs + '/>'

$start = $d.Range(0, 0)
$start.InsertXML($xml)
